$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.773136973381042
$ws.Range("B1").Value = 3.675705194473267
$ws.Range("C1").Value = 1.971169471740723
$ws.Range("D1").Value = 1.377047300338745
$ws.Range("E1").Value = 1.175451755523682
